$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark numeric-looking text cells as Text format so COM does not
# coerce them into actual numbers (they must stay text, matching the source).
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D9',
    'D10',
    'D11',
    'D12',
    'D14',
    'D16',
    'D19',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D33',
    'D35',
    'D36',
    'D37',
    'D38',
    'D42',
    'D44',
    'D46',
    'D47',
    'D49'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.333.66'
$ws.Range('D3').Value = '3.902.24'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '529.69'
$ws.Range('E5').Value = '  +9.73%  '
$ws.Range('D6').Value = '145.03'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').Value = '0.612'
$ws.Range('E7').Value = '  -1.56%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.718'
$ws.Range('E9').Value = '  -3.31%  '
$ws.Range('D10').Value = '0.174'
$ws.Range('E10').Value = '  -3.33%  '
$ws.Range('D11').Value = '0.0000337'
$ws.Range('E11').Value = '  -4.50%  '
$ws.Range('D12').Value = '42.21'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('D13').Value = '4.519.02'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').Value = '10.26'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '3.904.61'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '13.99'
$ws.Range('E16').Value = '  -1.74%  '
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('E18').Value = '  +6.57%  '
$ws.Range('D19').Value = '19.84'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = '69.243.62'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').Value = '425.58'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = '3.40'
$ws.Range('E22').Value = '  -5.32%  '
$ws.Range('D23').Value = '14.17'
$ws.Range('E23').Value = '  -4.04%  '
$ws.Range('D24').Value = '88.27'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('D25').Value = '4.03'
$ws.Range('E25').Value = '  +9.83%  '
$ws.Range('D26').Value = '11.42'
$ws.Range('E26').Value = '  -7.42%  '
$ws.Range('D27').Value = '10.60'
$ws.Range('E27').Value = '  -3.57%  '
$ws.Range('D28').Value = '36.43'
$ws.Range('E28').Value = '  -1.98%  '
$ws.Range('D29').Value = '690.03'
$ws.Range('E29').Value = '  -4.00%  '
$ws.Range('D30').Value = '13.19'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('E31').Value = '  -2.89%  '
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('D33').Value = '68.97'
$ws.Range('E33').Value = '  +11.66%  '
$ws.Range('D34').Value = '0.0₃0881'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').Value = '0.437'
$ws.Range('E35').Value = '  +9.79%  '
$ws.Range('D36').Value = '5.95'
$ws.Range('E36').Value = '  -1.73%  '
$ws.Range('D37').Value = '40.03'
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('D38').Value = '0.150'
$ws.Range('E38').Value = '  +2.43%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  +7.38%  '
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').Value = '3.23'
$ws.Range('E42').Value = '  +9.38%  '
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('D44').Value = '2.82'
$ws.Range('E44').Value = '  -5.73%  '
$ws.Range('E45').Value = '  +1.32%  '
$ws.Range('D46').Value = '0.000287'
$ws.Range('E46').Value = '  +14.58%  '
$ws.Range('D47').Value = '0.141'
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('E48').Value = '  +6.77%  '
$ws.Range('D49').Value = '146.00'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.746.36'
$ws.Range('E50').Value = '  +14.82%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0342'
$ws.Range('E51').Value = '  -2.82%  '

# Restore default styling on the cells we temporarily text-formatted,
# so no stray number-format / style id is left behind on them.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
